$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "XV-1701.9"
$ws.Range("C2").Value = "(NOTA 1)"
$ws.Range("D2").Value = "PHS2-PR27-1701-FL.01 Rev. 74I"
$ws.Range("F2").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G2").Value = "1120 kg/m3"
$ws.Range("H2").Value = "<50 cP"
$ws.Range("I2").Value = "1,2"
$ws.Range("J2").Value = "'2"
$ws.Range("L2").Value = "'65"
$ws.Range("M2").Value = "15 - 30"
$ws.Range("Q2").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R2").Value = "'15"
$ws.Range("S2").Value = "'30"

# Row 3
$ws.Range("A3").Value = "XV-1702.10"
$ws.Range("C3").Value = "BICICLOPIRONA EM PÓ"
$ws.Range("D3").Value = "PHS2-PR27-1701-FL.04 REV.00N"
$ws.Range("F3").Value = "GUILHOTINA AUTOMÁTICA"
$ws.Range("I3").Value = "-"
$ws.Range("J3").Value = "'2"
$ws.Range("L3").Value = "'65"
$ws.Range("N3").Value = "F"
$ws.Range("Q3").Value = "-"

# Row 4
$ws.Range("A4").Value = "XV-1702.9"
$ws.Range("C4").Value = "(NOTA 1)"
$ws.Range("D4").Value = "PHS2-PR27-1701-FL.01 Rev. 74I"
$ws.Range("F4").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G4").Value = "1120 kg/m3"
$ws.Range("H4").Value = "<50 cP"
$ws.Range("I4").Value = "1,2"
$ws.Range("J4").Value = "'2"
$ws.Range("L4").Value = "'65"
$ws.Range("M4").Value = "15 - 30"
$ws.Range("Q4").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R4").Value = "'15"
$ws.Range("S4").Value = "'30"

# Row 5
$ws.Range("A5").Value = "XV-1710.3"
$ws.Range("C5").Value = "(NOTA 1)"
$ws.Range("D5").Value = "PHS2-PR27-1701-FL.04 REV.00N"
$ws.Range("F5").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G5").Value = "1120 kg/m³"
$ws.Range("H5").Value = "<50 cP"
$ws.Range("I5").Value = "-"
$ws.Range("J5").Value = "'2"
$ws.Range("L5").Value = "'65"
$ws.Range("M5").Value = 15
$ws.Range("N5").Value = "F"
$ws.Range("Q5").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"

# Row 6
$ws.Range("A6").Value = "XV-1710.4"
$ws.Range("C6").Value = "(NOTA 1)"
$ws.Range("D6").Value = "PHS2-PR27-1701-FL.04 REV.00N"
$ws.Range("F6").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G6").Value = "1120 kg/m³"
$ws.Range("H6").Value = "<50 cP"
$ws.Range("I6").Value = "-"
$ws.Range("J6").Value = "'2"
$ws.Range("L6").Value = "'65"
$ws.Range("M6").Value = "15 - 30"
$ws.Range("N6").Value = "F"
$ws.Range("Q6").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R6").Value = "'15"
$ws.Range("S6").Value = "'30"

# Row 7
$ws.Range("A7").Value = "XV-1710.5"
$ws.Range("B7").Value = "3`""
$ws.Range("C7").Value = "(NOTA 1)"
$ws.Range("F7").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G7").Value = "1120 kg/m3"
$ws.Range("H7").Value = "<50 cP"
$ws.Range("I7").Value = "'0,984"
$ws.Range("J7").Value = "'2"
$ws.Range("L7").Value = "'65"
$ws.Range("M7").Value = "15 - 30"
$ws.Range("N7").Value = "F"
$ws.Range("Q7").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R7").Value = "'15"
$ws.Range("S7").Value = "'30"

# Row 8
$ws.Range("A8").Value = "XV-1710.6"
$ws.Range("C8").Value = "ÁGUA DE PROCESSO"
$ws.Range("D8").Value = "PHS2-PR27-1701-FL.04 Rev. 00N"
$ws.Range("F8").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G8").Value = "998 kg/m3"
$ws.Range("H8").Value = "1 cP"
$ws.Range("I8").Value = "'2"
$ws.Range("J8").Value = "'3"
$ws.Range("M8").Value = "1 - 6"
$ws.Range("N8").Value = "F"
$ws.Range("Q8").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R8").Value = "'1"
$ws.Range("S8").Value = "'6"

# Row 9
$ws.Range("A9").Value = "XV-1710.7"
$ws.Range("C9").Value = "ÁGUA DE PROCESSO"
$ws.Range("D9").Value = "PHS2-PR27-1701-FL.04 Rev. 00N"
$ws.Range("F9").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G9").Value = "998 kg/m3"
$ws.Range("H9").Value = "1 cP"
$ws.Range("I9").Value = "'2"
$ws.Range("J9").Value = "'3"
$ws.Range("M9").Value = "1 - 6"
$ws.Range("N9").Value = "F"
$ws.Range("Q9").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R9").Value = "'1"
$ws.Range("S9").Value = "'6"

# Row 10
$ws.Range("A10").Value = "XV-1710.8"
$ws.Range("C10").Value = "AR + VAPOR DA FASE LÍQUIDA DE REAÇÃO (NOTA 1)"
$ws.Range("D10").Value = "PHS2-PR27-1701-FL.04 Rev. 00N"
$ws.Range("F10").Value = "GUILHOTINA AUTOMÁTICA"
$ws.Range("G10").Value = "1,5 kg/m3"
$ws.Range("H10").Value = "0,0186 cP"
$ws.Range("I10").Value = "-"
$ws.Range("N10").Value = "A"
$ws.Range("Q10").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"

# Row 11
$ws.Range("A11").Value = "XV-1715.1"
$ws.Range("C11").Value = "(NOTA 1)"
$ws.Range("F11").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G11").Value = "1120 kg/m³"
$ws.Range("H11").Value = "<50 cP"
$ws.Range("I11").Value = "3,00"
$ws.Range("J11").Value = "6,00"
$ws.Range("M11").Value = "15 - 30"
$ws.Range("Q11").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R11").Value = "'15"
$ws.Range("S11").Value = "'30"

# Row 12
$ws.Range("A12").Value = "XV-1715.2"
$ws.Range("C12").Value = "(NOTA 1)"
$ws.Range("F12").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G12").Value = "1120 kg/m³"
$ws.Range("H12").Value = "<50 cP"
$ws.Range("I12").Value = "3,00"
$ws.Range("J12").Value = "6,00"
$ws.Range("M12").Value = "15 - 30"
$ws.Range("Q12").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R12").Value = "'15"
$ws.Range("S12").Value = "'30"

# Row 13
$ws.Range("A13").Value = "XV-1715.3"
$ws.Range("C13").Value = "(NOTA 1)"
$ws.Range("F13").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G13").Value = "1120 kg/m³"
$ws.Range("H13").Value = "<50 cP"
$ws.Range("I13").Value = "3,00"
$ws.Range("J13").Value = "6,00"
$ws.Range("K13").Value = "9 - 35"
$ws.Range("L13").Value = "'65"
$ws.Range("M13").Value = "15 - 30"
$ws.Range("Q13").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R13").Value = "'15"
$ws.Range("S13").Value = "'30"
$ws.Range("T13").Value = "'9"

# Row 14
$ws.Range("A14").Value = "XV-1715.4"
$ws.Range("B14").Value = "-"
$ws.Range("C14").Value = "(NOTA 1)"
$ws.Range("D14").Value = "PHS2-PR27-1701-FL.01 Rev. 74I"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G14").Value = "1120 kg/m³"
$ws.Range("H14").Value = "<50 cP"
$ws.Range("I14").Value = "3,00"
$ws.Range("J14").Value = "6,00"
$ws.Range("K14").Value = "9 - 35"
$ws.Range("L14").Value = "'65"
$ws.Range("M14").Value = "15 - 30"
$ws.Range("N14").Value = "-"
$ws.Range("O14").Value = "-"
$ws.Range("P14").Value = "-"
$ws.Range("Q14").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R14").Value = "'15"
$ws.Range("S14").Value = "'30"
$ws.Range("T14").Value = "'9"
$ws.Range("U14").Value = "'35"

# Row 15
$ws.Range("A15").Value = "XV-1715.5"
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "(NOTA 1)"
$ws.Range("D15").Value = "PHS2-PR27-1701-FL.01 Rev. 74I"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G15").Value = "1120 kg/m³"
$ws.Range("H15").Value = "<50 cP"
$ws.Range("I15").Value = "3,00"
$ws.Range("J15").Value = "6,00"
$ws.Range("K15").Value = "9 - 35"
$ws.Range("L15").Value = "'65"
$ws.Range("M15").Value = "15 - 30"
$ws.Range("N15").Value = "-"
$ws.Range("O15").Value = "-"
$ws.Range("P15").Value = "-"
$ws.Range("Q15").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R15").Value = "'15"
$ws.Range("S15").Value = "'30"
$ws.Range("T15").Value = "'9"
$ws.Range("U15").Value = "'35"

# Row 16
$ws.Range("A16").Value = "XV-1715.6"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "(NOTA 1)"
$ws.Range("D16").Value = "PHS2-PR27-1701-FL.01 Rev. 74I"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G16").Value = "1120 kg/m³"
$ws.Range("H16").Value = "<50 cP"
$ws.Range("I16").Value = "3,00"
$ws.Range("J16").Value = "6,00"
$ws.Range("K16").Value = "9 - 35"
$ws.Range("L16").Value = "'65"
$ws.Range("M16").Value = "15 - 30"
$ws.Range("N16").Value = "-"
$ws.Range("O16").Value = "-"
$ws.Range("P16").Value = "-"
$ws.Range("Q16").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R16").Value = "'15"
$ws.Range("S16").Value = "'30"
$ws.Range("T16").Value = "'9"
$ws.Range("U16").Value = "'35"

# Row 17
$ws.Range("A17").Value = "XV-1715.7"
$ws.Range("B17").Value = "-"
$ws.Range("C17").Value = "(NOTA 1)"
$ws.Range("D17").Value = "PHS2-PR27-1701-FL.01 Rev. 74I"
$ws.Range("E17").Value = "-"
$ws.Range("F17").Value = "ESFERA AUTOMÁTICA"
$ws.Range("G17").Value = "1120 kg/m³"
$ws.Range("H17").Value = "<50 cP"
$ws.Range("I17").Value = "3,00"
$ws.Range("J17").Value = "6,00"
$ws.Range("K17").Value = "9 - 35"
$ws.Range("L17").Value = "'65"
$ws.Range("M17").Value = "15 - 30"
$ws.Range("N17").Value = "-"
$ws.Range("O17").Value = "-"
$ws.Range("P17").Value = "-"
$ws.Range("Q17").Value = "(NOTA 1): CRESTIVO (A16003E); BIRK 25% (A25191A); BIRTEA; CRESTIVO DUO"
$ws.Range("R17").Value = "'15"
$ws.Range("S17").Value = "'30"
$ws.Range("T17").Value = "'9"
$ws.Range("U17").Value = "'35"

